# Auto-generated: apply cell value updates from diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "28.396.73"
Set-TextValue "E2" "  -0.05%  "
Set-TextValue "D3" "1.821.46"
Set-TextValue "E3" "  -0.11%  "
Set-TextValue "E4" "  +0.18%  "
Set-TextValue "D5" "315.72"
Set-TextValue "E5" "  +0.22%  "
Set-TextValue "E6" "  +0.14%  "
Set-TextValue "E7" "  +2.23%  "
Set-TextValue "D8" "0.3853"
Set-TextValue "E8" "  -1.76%  "
Set-TextValue "D9" "0.08043"
Set-TextValue "E9" "  +5.20%  "
Set-TextValue "B10" "OKB"
Set-TextValue "C10" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D10" "41.87"
Set-TextValue "E10" "  +0.63%  "
Set-TextValue "B11" "Polygon"
Set-TextValue "C11" "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-TextValue "D11" "1.114"
Set-TextValue "E11" "  +0.57%  "
Set-TextValue "D12" "6.392"
Set-TextValue "E12" "  +1.95%  "
Set-TextValue "B13" "Solana"
Set-TextValue "C13" "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D13" "20.91"
Set-TextValue "E13" "  -0.46%  "
Set-TextValue "B14" "BinanceUSD"
Set-TextValue "C14" "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue "D14" "1.004"
Set-TextValue "E14" "  +0.17%  "
Set-TextValue "D15" "7.429"
Set-TextValue "E15" "  -1.05%  "
Set-TextValue "D16" "1.818.89"
Set-TextValue "E16" "  +0.02%  "
Set-TextValue "D17" "94.42"
Set-TextValue "E17" "  +1.27%  "
Set-TextValue "D18" "0.00001106"
Set-TextValue "E18" "  +1.04%  "
Set-TextValue "E19" "  -0.70%  "
Set-TextValue "D20" "17.62"
Set-TextValue "E20" "  -0.15%  "
Set-TextValue "E21" "  +0.21%  "
Set-TextValue "D22" "6.017"
Set-TextValue "E22" "  -2.05%  "
Set-TextValue "D23" "28.446.94"
Set-TextValue "E23" "  +0.05%  "
Set-TextValue "E24" "  +1.67%  "
Set-TextValue "D25" "2.250"
Set-TextValue "E25" "  -0.34%  "
Set-TextValue "D26" "158.95"
Set-TextValue "E26" "  +1.75%  "
Set-TextValue "D27" "20.85"
Set-TextValue "E27" "  +0.62%  "
Set-TextValue "D28" "2.029.63"
Set-TextValue "E28" "  -0.04%  "
Set-TextValue "D29" "2.410"
Set-TextValue "E29" "  +1.20%  "
Set-TextValue "E30" "  +0.22%  "
Set-TextValue "E31" "  +2.01%  "
Set-TextValue "D32" "1.078"
Set-TextValue "E32" "  -2.69%  "
Set-TextValue "D33" "5.674"
Set-TextValue "E33" "  +0.63%  "
Set-TextValue "D34" "3.680"
Set-TextValue "E34" "  +0.56%  "
Set-TextValue "D35" "0.07298"
Set-TextValue "E35" "  +3.62%  "
Set-TextValue "E36" "  +8.76%  "
Set-TextValue "D37" "0.2204"
Set-TextValue "E37" "  -0.11%  "
Set-TextValue "D38" "0.02344"
Set-TextValue "E38" "  +1.10%  "
Set-TextValue "D39" "8.790"
Set-TextValue "E39" "  -0.44%  "
Set-TextValue "D40" "5.107"
Set-TextValue "E40" "  -0.95%  "
Set-TextValue "E41" "  +1.13%  "
Set-TextValue "D42" "1.183"
Set-TextValue "E42" "  +1.02%  "
Set-TextValue "E43" "  -0.35%  "
Set-TextValue "B44" "EnergySwap"
Set-TextValue "C44" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D44" "13.51"
Set-TextValue "E44" "  +1.18%  "
Set-TextValue "B45" "Decentraland"
Set-TextValue "C45" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D45" "0.6143"
Set-TextValue "E45" "  +4.44%  "
Set-TextValue "D46" "3.804"
Set-TextValue "E46" "  +2.62%  "
Set-TextValue "D47" "127.02"
Set-TextValue "E47" "  +1.48%  "
Set-TextValue "D48" "1.218"
Set-TextValue "E48" "  +1.99%  "
Set-TextValue "D49" "1.968"
Set-TextValue "E49" "  -0.39%  "
Set-TextValue "D50" "0.06903"
Set-TextValue "E50" "  -0.07%  "
Set-TextValue "D51" "73.95"
Set-TextValue "E51" "  -0.02%  "
